# Updated cryptos list on Tue Aug 22 13:25:10 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures on Sheet1
# for rows 2-51. Price cells hold numeric-looking text (e.g. "1.666.27"), so
# each is written with a leading apostrophe to keep it stored as text instead
# of being auto-converted to a number, matching the original inline-string data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "'1.666.27"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'208.59"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("D6").Value = "'0.5210"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'0.2600"
$ws.Range("E8").Value = "  -2.76%  "
$ws.Range("D9").Value = "'0.06331"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "'21.01"
$ws.Range("E10").Value = "  -0.93%  "
$ws.Range("D11").Value = "'0.07544"
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("D12").Value = "'1.679.59"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "'4.413"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").Value = "'0.5388"
$ws.Range("E14").Value = "  -4.56%  "
$ws.Range("D15").Value = "'0.0₅7994"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("D16").Value = "'66.12"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "'26.162.51"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").Value = "'4.713"
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("D20").Value = "'187.43"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "'10.23"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").Value = "'6.213"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "'1.005"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'149.10"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").Value = "'7.420"
$ws.Range("E26").Value = "  -2.16%  "
$ws.Range("D27").Value = "'15.70"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").Value = "'0.06190"
$ws.Range("E28").Value = "  -3.74%  "
$ws.Range("E29").Value = "  +2.43%  "
$ws.Range("D30").Value = "'1.271"
$ws.Range("E30").Value = "  -0.76%  "
$ws.Range("D31").Value = "'3.480"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").Value = "'3.397"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("D33").Value = "'1.636"
$ws.Range("E33").Value = "  -0.62%  "
$ws.Range("D34").Value = "'0.9922"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").Value = "'2.394"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("E36").Value = "  +1.82%  "
$ws.Range("D37").Value = "'0.5896"
$ws.Range("E37").Value = "  -3.06%  "
$ws.Range("D38").Value = "'1.107.95"
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("D39").Value = "'6.025"
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("D40").Value = "'0.01597"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'0.8539"
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").Value = "'99.98"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "'1.818.89"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").Value = "'0.0₈112"
$ws.Range("E45").Value = "  +1.93%  "
$ws.Range("D46").Value = "'55.39"
$ws.Range("E46").Value = "  -2.51%  "
$ws.Range("D47").Value = "'1.002"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").Value = "'8.041"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D49").Value = "'0.05256"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "'0.4260"
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").Value = "'5.879"
$ws.Range("E51").Value = "  -1.04%  "
